$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - F column ("想去人数" / want-to-go count) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1307
$ws1.Range("F3").Value = 1188
$ws1.Range("F11").Value = 2365
$ws1.Range("F13").Value = 1365
$ws1.Range("F24").Value = 4706
$ws1.Range("F26").Value = 204
$ws1.Range("F29").Value = 127
$ws1.Range("F31").Value = 88
$ws1.Range("F33").Value = 683
$ws1.Range("F39").Value = 977

# Sheet "全部类型" (all types) - same underlying events, F column updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1307
$ws4.Range("F5").Value = 1188
$ws4.Range("F18").Value = 2365
$ws4.Range("F20").Value = 1365
$ws4.Range("F30").Value = 4706
$ws4.Range("F34").Value = 127
$ws4.Range("F36").Value = 88
$ws4.Range("F38").Value = 683
$ws4.Range("F42").Value = 977
